# "Generate Report for Handoff" - refresh localization status report.
#
# The four files that were at "low" priority moved to "ht" (high) priority,
# and the zh-cn handoff pass re-generated its xliff handoff batch, bumping
# the "Latest Handoff Datetime" for those same four rows as well as the
# "Latest HO Xliff Generate Date" on the Overview sheet for the one file
# that is still "In Translation".

$wb = $excel.ActiveWorkbook

# --- Overview sheet: bump the Latest HO Xliff Generate Date shared by
#     rows 4-7 (these four rows all shared the same date string).
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($row in 4..7) {
    $wsOverview.Cells.Item($row, 7).Value = "2016-10-26 08:58:19"
}

# --- zh-cn sheet: Priority low -> ht, and Latest Handoff Datetime bump
#     for rows 4-7 (the four files that were "low" priority).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($row in 4..7) {
    $wsZhCn.Cells.Item($row, 5).Value = "ht"
    $wsZhCn.Cells.Item($row, 8).Value = "2016-10-26 08:58:06"
}

# --- de-de sheet: Priority low -> ht, and Latest Handoff Datetime bump
#     (shares the same underlying date string as the Overview sheet)
#     for rows 4-7.
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($row in 4..7) {
    $wsDeDe.Cells.Item($row, 5).Value = "ht"
    $wsDeDe.Cells.Item($row, 8).Value = "2016-10-26 08:58:19"
}
